# "10Th - MB for single stock and added new group"
#
# The weekly MarketBeat-rank watch sheet rolls forward by one column each
# update: the three oldest date columns (B:D -> Jun_17/Jun_15/Jun_13) slide
# right to make room for three new/refreshed columns on the left, and two
# new rows are appended for a newly-tracked benchmark/analyst group.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing B:D columns (and everything to their right) three
# columns to the right, opening up B:D for the new week's data.
$ws.Columns("B:D").Insert()

# Header row: the two freshly-inserted "inner" columns both restate the
# prior business day (Jun_26); B1 gets the new, most-recent date (Jun_27).
$ws.Range("B1").Value = "Jun_27"
$ws.Range("C1").Value = "Jun_26"
$ws.Range("D1").Value = "Jun_26"

# Existing rating rows (2-27): fill the three newly-opened cells with the
# same "UN" (unchanged) marker used across the rest of the row.
$ws.Range("B2:D27").Value = "UN"

# New rows for the newly-tracked group.
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28:D28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29:D29").Value = "UN"
